# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap display order of "Santa Lucia" / "Timor Oriental" ---
$ws.Range("A202").Value = "Santa Lucia"
$ws.Range("A203").Value = "Timor Oriental"

# --- Swap display order of "Islas Malvinas" / "Montserrat" (and their stats) ---
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1

# --- Update "last updated" timestamp text ---
$ws.Range("A1").Value = "Datos actualizados a 10 de Agosto de 2020 a las 12:34"

# --- Row 4: Estados Unidos ---
$ws.Range("B4").Value = 5199524
$ws.Range("C4").Value = 80
$ws.Range("D4").Value = 2664849
$ws.Range("E4").Value = 2369058

# --- Row 36: Oman ---
$ws.Range("B36").Value = 81787
$ws.Range("C36").Value = 207
$ws.Range("D36").Value = 76124
$ws.Range("E36").Value = 5142
$ws.Range("G36").Value = 8
$ws.Range("H36").Value = 521

# --- Row 59: Argelia ---
$ws.Range("B59").Value = 35214
$ws.Range("E59").Value = 9406

# --- Row 87: Consejo Danes para los Refugiados ---
$ws.Range("B87").Value = 9489
$ws.Range("C87").Value = 35
$ws.Range("D87").Value = 8363
$ws.Range("E87").Value = 902

# --- Row 88: Malasia ---
$ws.Range("B88").Value = 9094
$ws.Range("C88").Value = 11
$ws.Range("D88").Value = 8803
$ws.Range("E88").Value = 166

# --- Row 95: Finlandia ---
$ws.Range("B95").Value = 7601
$ws.Range("C95").Value = 17
$ws.Range("E95").Value = 290

# --- Row 174: Islas Feroe ---
$ws.Range("B174").Value = 306
$ws.Range("C174").Value = 3
$ws.Range("E174").Value = 91

# --- Row 186: Belice ---
$ws.Range("B186").Value = 154
$ws.Range("C186").Value = 1
$ws.Range("E186").Value = 120
